# "voided product group by customer sku"
# Replace the Catalog No. (SKU) values for the voided rows with placeholder
# values, reset the Qty on those same rows to 5, and mark the Catalog No.
# column so it wraps its text (matches the new cellXfs style introduced by
# this edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Catalog No. (column G) value swaps -----------------------------------
$ws.Range("G2").Value = "iphone123"
$ws.Range("G3").Value = "iphone456"
$ws.Range("G4").Value = "iphone789"
$ws.Range("G5").Value = "MRI21"

# --- Qty (column J) reset to 5 for the voided rows -------------------------
$ws.Range("J2").Value = 5
$ws.Range("J3").Value = 5
$ws.Range("J4").Value = 5

# --- Formatting: wrap text on the Catalog No. column for these rows --------
$fmtRange = $ws.Range("G2:G5")
$fmtRange.Font.Name = "ARIAL"
$fmtRange.Font.Color = 0
$fmtRange.Font.Size = 10
$fmtRange.VerticalAlignment = -4107
$fmtRange.WrapText = $true

# --- Cursor position --------------------------------------------------------
$ws.Range("G5").Select() | Out-Null
